# Portal Check Added and Asana Updates
#
# Swap the SKU/UNITS/COST details between data row 3 and data row 5
# (columns I, J, L, M, N, O, P) on the active sheet. The affected cells
# hold text-like values (e.g. "100288", "13.78") that must stay stored as
# text/shared-strings rather than being re-interpreted as numbers, so the
# swap is performed with Copy / PasteSpecial(values) through a scratch
# cell instead of a plain Value assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

$rowA = 3
$rowB = 5
$cols = @("I", "J", "L", "M", "N", "O", "P")

$scratch = $ws.Range("Z1000")

foreach ($col in $cols) {
    $cellA = $ws.Range($col + $rowA)
    $cellB = $ws.Range($col + $rowB)

    $cellA.Copy()
    $scratch.PasteSpecial($xlPasteValues)

    $cellB.Copy()
    $cellA.PasteSpecial($xlPasteValues)

    $scratch.Copy()
    $cellB.PasteSpecial($xlPasteValues)
}

$scratch.ClearContents()
$excel.CutCopyMode = 0
